# Insert two new weekly-price rows for "Betarraga" (Agrícola del Norte S.A. de
# Arica) at the top of the data block (row 64), pushing all later rows down by
# two. This mirrors Excel's normal "insert rows" behaviour: formatting/styles
# of the following rows travel with them, and the sheet's used range grows by
# two rows (A1:R188 -> A1:R190).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64 downwards (row 64 becomes 66, ..., old 187/188 become 189/190)
$ws.Rows("64:65").Insert()

# --- New row 64 (Primera) ---
$ws.Cells.Item(64, 1).Value  = 1
$ws.Cells.Item(64, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value  = 44498
$ws.Cells.Item(64, 5).Value  = 15
$ws.Cells.Item(64, 6).Value  = 100114014
$ws.Cells.Item(64, 7).Value  = "Betarraga"
$ws.Cells.Item(64, 8).Value  = "Sin especificar"
$ws.Cells.Item(64, 9).Value  = "Primera"
$ws.Cells.Item(64, 10).Value = 1200
$ws.Cells.Item(64, 11).Value = 350
$ws.Cells.Item(64, 12).Value = 400
$ws.Cells.Item(64, 13).Value = 375
$ws.Cells.Item(64, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(64, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value = 94
$ws.Cells.Item(64, 17).Value = 4
$ws.Cells.Item(64, 18).Value = "Hortaliza"

# --- New row 65 (Segunda) ---
$ws.Cells.Item(65, 1).Value  = 1
$ws.Cells.Item(65, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(65, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(65, 4).Value  = 44498
$ws.Cells.Item(65, 5).Value  = 15
$ws.Cells.Item(65, 6).Value  = 100114014
$ws.Cells.Item(65, 7).Value  = "Betarraga"
$ws.Cells.Item(65, 8).Value  = "Sin especificar"
$ws.Cells.Item(65, 9).Value  = "Segunda"
$ws.Cells.Item(65, 10).Value = 1200
$ws.Cells.Item(65, 11).Value = 350
$ws.Cells.Item(65, 12).Value = 400
$ws.Cells.Item(65, 13).Value = 375
$ws.Cells.Item(65, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(65, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(65, 16).Value = 75
$ws.Cells.Item(65, 17).Value = 5
$ws.Cells.Item(65, 18).Value = "Hortaliza"
